$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (columns B,C,D,E,F,I,J,K,L,M,N), keyed by row number.
# Row 2 corresponds to bus index 0 (column A = 0) through row 25 = bus index 23.
# Column H (bus index 6) has no data in this sheet (left blank, as in the source).
$rowData = @{
    2 = @{ "B"=1.02; "C"=1.069207786466924; "D"=1.072378388838149; "E"=1.07150096090519; "F"=1.082863914496079; "I"=1.054163625004125; "J"=1.074143681852799; "K"=1.075073061932072; "L"=1.074197967000746; "M"=1.08553102344964; "N"=1.075669089087033 }
    3 = @{ "B"=1.02; "C"=1.070455277658043; "D"=1.073396408285127; "E"=1.07264728027479; "F"=1.084051031842779; "I"=1.054553841754828; "J"=1.075046830601396; "K"=1.075907318026839; "L"=1.075160037499017; "M"=1.086535962520786; "N"=1.076573520410445 }
    4 = @{ "B"=1.02; "C"=1.07126211954777; "D"=1.074054817702858; "E"=1.073389450779015; "F"=1.084819407441205; "I"=1.054804946869104; "J"=1.075630310466; "K"=1.076446207677332; "L"=1.075782398852134; "M"=1.087185874465074; "N"=1.077157828883383 }
    5 = @{ "B"=1.02; "C"=1.071601229949245; "D"=1.074331538206241; "E"=1.073701562120818; "F"=1.085142489216777; "I"=1.05491017902914; "J"=1.075875387060044; "K"=1.076672535427474; "L"=1.076044001638469; "M"=1.087459014650369; "N"=1.077403253514301 }
    6 = @{ "B"=1.02; "C"=1.071658163104293; "D"=1.074377996442016; "E"=1.073753973051244; "F"=1.085196739429747; "I"=1.054927828496318; "J"=1.075916523734888; "K"=1.076710523899388; "L"=1.076087923702388; "M"=1.087504871268857; "N"=1.077444448607941 }
    7 = @{ "B"=1.02; "C"=1.071266651092375; "D"=1.074058515548329; "E"=1.07339362082053; "F"=1.084823724251382; "I"=1.054806354291811; "J"=1.075633586047658; "K"=1.076449232748561; "L"=1.075785894549145; "M"=1.087189524501438; "N"=1.077161109116743 }
    8 = @{ "B"=1.02; "C"=1.069629459574768; "D"=1.072722499644674; "E"=1.071888277499122; "F"=1.083265059254316; "I"=1.054295789021958; "J"=1.074449096058563; "K"=1.075355195901342; "L"=1.074523137496834; "M"=1.085870720675424; "N"=1.075974937016 }
    9 = @{ "B"=1.02; "C"=1.066741597065904; "D"=1.070365788146631; "E"=1.069238878889151; "F"=1.08052020514547; "I"=1.053385424327093; "J"=1.072354782498356; "K"=1.073420190263292; "L"=1.072296707920741; "M"=1.083544079035063; "N"=1.073877649290316 }
    10 = @{ "B"=1.02; "C"=1.064814247612831; "D"=1.068792905176694; "E"=1.0674747007018; "F"=1.078691369281623; "I"=1.052771290599237; "J"=1.070953717330613; "K"=1.072125289593025; "L"=1.070811488676055; "M"=1.081991064184378; "N"=1.07247459444933 }
    11 = @{ "B"=1.02; "C"=1.063979151514958; "D"=1.06811139903261; "E"=1.066711270864525; "F"=1.077899696537591; "I"=1.052503640960904; "J"=1.070345868431243; "K"=1.071563404622747; "L"=1.070168135813366; "M"=1.081318118073078; "N"=1.071865882334803 }
    12 = @{ "B"=1.02; "C"=1.063668875450595; "D"=1.067858189997397; "E"=1.06642776828423; "F"=1.077605666078925; "I"=1.052403963816418; "J"=1.07011990738589; "K"=1.071354516003255; "L"=1.069929127980798; "M"=1.081068082094609; "N"=1.071639600398846 }
    13 = @{ "B"=1.02; "C"=1.06373543457614; "D"=1.067912507266453; "E"=1.066488577440424; "F"=1.077668735189061; "I"=1.052425356688954; "J"=1.070168384922553; "K"=1.071399331470395; "L"=1.069980397711444; "M"=1.081121719028262; "N"=1.071688146779173 }
    14 = @{ "B"=1.02; "C"=1.063953505732962; "D"=1.06809047007707; "E"=1.06668783504877; "F"=1.077875391244786; "I"=1.052495406929201; "J"=1.070327194080978; "K"=1.071546141495653; "L"=1.07014838013336; "M"=1.081297451536765; "N"=1.071847181464818 }
    15 = @{ "B"=1.02; "C"=1.064087855326914; "D"=1.068200109871449; "E"=1.066810613369251; "F"=1.078002723054108; "I"=1.052538532694969; "J"=1.070425017883403; "K"=1.071636572161734; "L"=1.070251874578846; "M"=1.081405716302555; "N"=1.071945144188263 }
    16 = @{ "B"=1.02; "C"=1.064869658543329; "D"=1.068838125143162; "E"=1.067525376820167; "F"=1.078743914587152; "I"=1.052789017190397; "J"=1.070994033273816; "K"=1.072162554958907; "L"=1.070854180644295; "M"=1.082035715138577; "N"=1.072514967645798 }
    17 = @{ "B"=1.02; "C"=1.065359916139815; "D"=1.069238217031639; "E"=1.067973854077841; "F"=1.079208903494984; "I"=1.052945676856275; "J"=1.071350644787039; "K"=1.072492171852921; "L"=1.071231925241523; "M"=1.082430766974451; "N"=1.072872085588273 }
    18 = @{ "B"=1.02; "C"=1.065645823093779; "D"=1.069471541763467; "E"=1.068235488818398; "F"=1.079480145514336; "I"=1.053036887346255; "J"=1.071558536673725; "K"=1.072684317637351; "L"=1.071452233912052; "M"=1.082661147784605; "N"=1.073080272705284 }
    19 = @{ "B"=1.02; "C"=1.065743301242422; "D"=1.069551092374395; "E"=1.068324707273622; "F"=1.079572635794336; "I"=1.05306795956129; "J"=1.071629403225866; "K"=1.072749815065808; "L"=1.071527349554866; "M"=1.082739693893265; "N"=1.073151239896058 }
    20 = @{ "B"=1.02; "C"=1.065307321541102; "D"=1.069195295325267; "E"=1.067925732013513; "F"=1.079159012339129; "I"=1.052928885979611; "J"=1.071312395508623; "K"=1.0724568189014; "L"=1.071191399248711; "M"=1.082388386447128; "N"=1.072833781991495 }
    21 = @{ "B"=1.02; "C"=1.063889291559977; "D"=1.068038066335805; "E"=1.066629156749853; "F"=1.077814535316905; "I"=1.052474786068001; "J"=1.070280433695448; "K"=1.071502914570765; "L"=1.070098914555442; "M"=1.08124570474437; "N"=1.071800354674178 }
    22 = @{ "B"=1.02; "C"=1.062997231285331; "D"=1.067310080318416; "E"=1.065814348349443; "F"=1.076969393776884; "I"=1.052187769832376; "J"=1.069630562365976; "K"=1.070902117606251; "L"=1.069411805893983; "M"=1.080526827428072; "N"=1.071149560452876 }
    23 = @{ "B"=1.02; "C"=1.063470176594487; "D"=1.067696037000993; "E"=1.066246256293948; "F"=1.07741740240761; "I"=1.052340065579564; "J"=1.06997517027275; "K"=1.071220710500067; "L"=1.069776076536865; "M"=1.08090795901477; "N"=1.071494657742405 }
    24 = @{ "B"=1.02; "C"=1.065331086937029; "D"=1.069214689926987; "E"=1.067947476162379; "F"=1.079181555937476; "I"=1.052936473568325; "J"=1.071329679062077; "K"=1.072472793731263; "L"=1.07120971127271; "M"=1.082407536526126; "N"=1.072851090089577 }
    25 = @{ "B"=1.02; "C"=1.06748853974025; "D"=1.07097535614834; "E"=1.069923438863216; "F"=1.081229621459584; "I"=1.053622045681393; "J"=1.072897062169252; "K"=1.073921293261803; "L"=1.07287245295634; "M"=1.074420699060546 }
}

foreach ($row in $rowData.Keys) {
    $cellValues = $rowData[$row]
    foreach ($col in $cellValues.Keys) {
        $ws.Range("$col$row").Value = $cellValues[$col]
    }
}
